# Re-running university responses / police actions analyses.
# The underlying event coding was re-run, which changed two linked numbers:
#   - summary_counts!B13 ("events with neither university response nor
#     police coding") went from 4942 -> 4984 (+42 events).
#   - response_action_counts!E (the "percentage" column) is each row's
#     count (column D) divided by the total number of coded events,
#     which is summary_counts!B13 + summary_counts!B14. That denominator
#     moves from 5365 to 5407, so every percentage in E2:E73 is
#     recomputed against the new denominator.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("summary_counts")
$responses = $wb.Worksheets.Item("response_action_counts")

# Update the "neither university response nor police coding" count.
$summary.Range("B13").Value2 = 4984

# New denominator used for the percentage column: events with neither +
# total number of protests (the same quantity the original percentages
# were computed against, just re-derived from the refreshed counts).
$denominator = $summary.Range("B13").Value2 + $summary.Range("B14").Value2

# Recompute column E (percentage) for every data row, based on the
# existing column D (count) values, using the refreshed denominator.
# The source analysis keeps ~16 significant decimal digits of precision
# (matching how the numbers were originally produced), so round-trip the
# ratio through a 16-significant-digit string before storing it - this
# keeps the stored value identical to what the upstream tooling emits.
$lastRow = $responses.Cells.Item($responses.Rows.Count, "D").End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $countCell = $responses.Cells.Item($r, 4)
    $pctCell = $responses.Cells.Item($r, 5)
    if ($countCell.Value2 -ne $null -and $pctCell.Value2 -ne $null) {
        $ratio = $countCell.Value2 / $denominator
        $pctCell.Value2 = $ratio.ToString("G16")
    }
}
